$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "09406569642"
$ws.Range("B8").Value = "000000"
$ws.Range("C8").Value = "7d884018af47058d31"
$ws.Range("D8").Value = "dbwew6HPhkY:APA91bH6tZ_jD97MM47kzxjjhLjsWN17CoUFSwYrAZe3-_XHt0plqPfZhk8T2FpPecdp6zP7vGQACVN9l9IyhxAcFs6yecpJ03vNI_fvSTGbHzmoG5KFnxA5TsZmda5u-oyK770AZyRU"
$ws.Range("E8").Value = "5cf9d8f7ebfc2523d81a1d9c"
$ws.Range("F8").Formula = "=A8"
$ws.Range("G8").Value = "customer"
$ws.Range("H8").Value = "09406569642"
$ws.Range("I8").Value = "customer"
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = "MMK"
$ws.Range("L8").Value = "Abc"
$ws.Range("M8").Value = "FO"
$ws.Range("N8").Formula = "=B8"
$ws.Range("O8").Value = 1234
$ws.Range("P8").Value = 9081
$ws.Range("Q8").Value = "Receiver need to be another account"

$ws.Range("F8:H8").Style = "Bad"

$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("P9").Select()
